$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "youth basketball leggings with knee pads",
    "youth basketball compression pants with knee pads",
    "softball sliding knee pads",
    "baseball sliding knee pads",
    "hex pads basketball knee youth",
    "youth basketball knee pads",
    "sliding knee pads softball",
    "youth black knee pads for basketball",
    "crossfit leggings",
    "knee pads for basketball",
    "basketball pads",
    "basketball knee",
    "basketball leggings for boys",
    "basketball leg pads tights",
    "compression basketball pants men",
    "crossfit knee sleeves",
    "padded leg sleeves for basketball youth",
    "hex knee pad",
    "compression pants men capri",
    "knee sleeves for crossfit",
    "basketball leg compression pants",
    "compression leggings for men",
    "snowboarding knee pads men",
    "mens tights and leggings",
    "adult football pants with pads",
    "volleyball knee pads xxl",
    "compression basketball pants boys",
    "basketball knee compression sleeve",
    "padded knee compression sleeve",
    "knee pads basketball girls",
    "knee pads for work pants",
    "knee brace wrestling",
    "crossfit knee support",
    "men compression leggings",
    "basketball training gear",
    "sports knee pad",
    "hex padded compression leg sleeve",
    "knee sleeve weightlifting men",
    "hex pads",
    "cycling capris padded",
    "knee pads volleyball men",
    "boy compression pants basketball",
    "crossfit knee compression sleeve",
    "crossfit knee sleeve",
    "mens tights",
    "protector de rodillas basketball",
    "workout pads",
    "mens athletic tights",
    "athletic capri",
    "mens basketball leggings",
    "mens sliding pants baseball",
    "padded compression pants",
    "sliding pants baseball",
    "tights capri men",
    "tights with knee pads basketball",
    "workout knee pads for men",
    "workout sliding pads",
    "wrestling knee pad sleeve",
    "youth knee pads basketball",
    "youth wrestling knee pads",
    "leggings knee pads",
    "xl knee pads",
    "gel knee pads wrestling",
    "men s leggings compression",
    "capri tights for men",
    "padded knee pads",
    "football knee pads youth",
    "knee compression sleeve cycling",
    "volleyball knee pads for girls",
    "youth football leggings",
    "youth basketball leggings boys",
    "wrestling knee sleeve",
    "compression sleeve youth",
    "workout capri pants",
    "workout tights",
    "basketball pants",
    "compression leggings men",
    "youth knee pads",
    "basketball leggins with knee pads",
    "knee pads black mens",
    "basketball knee pads white",
    "knee pads white basketball",
    "black workout capris",
    "black workout leggings capri",
    "workout pads for knees",
    "knee pads for basketball blue",
    "knee tights for men",
    "mcdavid padded leg sleeve",
    "teen leggings",
    "workout pants for men",
    "mens training pants",
    "basketball knee pads women",
    "legging for men workout",
    "mens athletic leggings",
    "mens capri",
    "workout support",
    "athletic legging",
    "basketball knee pads boys youth",
    "basketball spandex knee pads",
    "training gear six pad"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
